$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Swap China (row5) and Colombia (row6) rows, and move the Korea, Rep. row
# (row13) down past Singapore so the table reflects the reordered country
# list exactly like the authored edit did.

$ws.Cells.Item(5,1).Value = "Colombia"
$ws.Cells.Item(5,2).Value = 142.47
$ws.Cells.Item(5,3).Value = 147.94999999999999
$ws.Cells.Item(5,4).Value = 136.99
$ws.Cells.Item(5,5).Value = 136.99
$ws.Cells.Item(5,6).Value = 134.25

$ws.Cells.Item(6,1).Value = "China"
$ws.Cells.Item(6,2).Value = 16.440000000000001
$ws.Cells.Item(6,3).Value = 19.18
$ws.Cells.Item(6,4).Value = 16.440000000000001
$ws.Cells.Item(6,5).Value = 19.18
$ws.Cells.Item(6,6).Value = 21.92

$ws.Cells.Item(13,1).Value = "Lebanon"
$ws.Cells.Item(13,2).Value = 90.41
$ws.Cells.Item(13,3).Value = 104.11
$ws.Cells.Item(13,4).Value = 90.41
$ws.Cells.Item(13,5).Value = 82.19
$ws.Cells.Item(13,6).Value = 93.15

$ws.Cells.Item(14,1).Value = "Malaysia"
$ws.Cells.Item(14,2).Value = 95.89
$ws.Cells.Item(14,3).Value = 98.63
$ws.Cells.Item(14,4).Value = 120.55
$ws.Cells.Item(14,5).Value = 115.07
$ws.Cells.Item(14,6).Value = 112.33

$ws.Cells.Item(15,1).Value = "Mexico"
$ws.Cells.Item(15,2).Value = 128.77000000000001
$ws.Cells.Item(15,3).Value = 134.25
$ws.Cells.Item(15,4).Value = 136.99
$ws.Cells.Item(15,5).Value = 131.51
$ws.Cells.Item(15,6).Value = 131.51

$ws.Cells.Item(16,1).Value = "Netherlands"
$ws.Cells.Item(16,2).Value = 98.63
$ws.Cells.Item(16,3).Value = 142.47
$ws.Cells.Item(16,4).Value = 134.25
$ws.Cells.Item(16,5).Value = 123.29
$ws.Cells.Item(16,6).Value = 142.47

$ws.Cells.Item(17,1).Value = "New Zealand"
$ws.Cells.Item(17,2).Value = 128.77000000000001
$ws.Cells.Item(17,3).Value = 142.47
$ws.Cells.Item(17,4).Value = 145.21
$ws.Cells.Item(17,5).Value = 158.9
$ws.Cells.Item(17,6).Value = 164.38

$ws.Cells.Item(18,1).Value = "Pakistan"
$ws.Cells.Item(18,2).Value = 82.19
$ws.Cells.Item(18,3).Value = 82.19
$ws.Cells.Item(18,4).Value = 68.489999999999995
$ws.Cells.Item(18,5).Value = 71.23
$ws.Cells.Item(18,6).Value = 73.97

$ws.Cells.Item(19,1).Value = "Philippines"
$ws.Cells.Item(19,2).Value = 65.75
$ws.Cells.Item(19,3).Value = 71.23
$ws.Cells.Item(19,4).Value = 79.45
$ws.Cells.Item(19,5).Value = 79.45
$ws.Cells.Item(19,6).Value = 76.709999999999994

$ws.Cells.Item(20,1).Value = "Poland"
$ws.Cells.Item(20,2).Value = 134.25
$ws.Cells.Item(20,3).Value = 139.72999999999999
$ws.Cells.Item(20,4).Value = 115.07
$ws.Cells.Item(20,5).Value = 117.81
$ws.Cells.Item(20,6).Value = 123.29

$ws.Cells.Item(21,1).Value = "Romania"
$ws.Cells.Item(21,2).Value = 79.45
$ws.Cells.Item(21,3).Value = 68.489999999999995
$ws.Cells.Item(21,4).Value = 73.97
$ws.Cells.Item(21,5).Value = 60.27
$ws.Cells.Item(21,6).Value = 71.23

$ws.Cells.Item(22,1).Value = "Russia"
$ws.Cells.Item(22,2).Value = ""
$ws.Cells.Item(22,3).Value = ""
$ws.Cells.Item(22,4).Value = 95.89
$ws.Cells.Item(22,5).Value = 109.59
$ws.Cells.Item(22,6).Value = 120.55

$ws.Cells.Item(23,1).Value = "Singapore"
$ws.Cells.Item(23,2).Value = ""
$ws.Cells.Item(23,3).Value = ""
$ws.Cells.Item(23,4).Value = ""
$ws.Cells.Item(23,5).Value = ""
$ws.Cells.Item(23,6).Value = ""

$ws.Cells.Item(24,1).Value = "Korea, Rep."
$ws.Cells.Item(24,2).Value = 54.8
$ws.Cells.Item(24,3).Value = 84.93
$ws.Cells.Item(24,4).Value = 82.19
$ws.Cells.Item(24,5).Value = 87.67
$ws.Cells.Item(24,6).Value = 98.63


# Restore the selection to B1:H30 with B1 active, matching the saved view.
$ws.Range("B1:H30").Select()
